# resultados.xlsx - "adicionando tabelas de resultados"
#
# 1. Fix a stray trailing-space typo in the Berea500 "Poros"/image-id
#    label (the sample that used to read "I31 " now reads "I31").
# 2. Re-sort the Berea500 results table (A1:E8, header included) in
#    ascending order by the image-id column (A) - this is what turns
#    the I32/I318/I320/I311/I312/I310/I31 ordering into the alphabetic
#    I31/I310/I311/I312/I318/I32/I320 ordering.
# 3. Apply a 6-decimal numeric format to the "Porosidade (%)" column on
#    both the Berea200 and Berea500 tables.
# 4. Leave the workbook with the Berea500 tab active.

$wb = $excel.ActiveWorkbook

$berea200 = $wb.Worksheets.Item("Berea200")
$berea500 = $wb.Worksheets.Item("Berea500")

# --- fix the "I31 " -> "I31" label before we sort on it ---------------
$berea500.Range("A8").Value = "I31"

# --- sort the Berea500 table by column A (ascending), header intact ---
$sortRange = $berea500.Range("A1:E8")
$sortKey = $berea500.Range("A2:A8")
$sortRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1)

# --- number formats: Porosidade column gets 6 decimal places ----------
$berea200.Range("B2:B5").NumberFormat = "0.000000"
$berea500.Range("B2:B8").NumberFormat = "0.000000"

# --- make Berea500 the active tab --------------------------------------
$berea500.Activate()
